# Add 2022-Q4 data
# ------------------------------------------------------------------
# Before: sheets = [ 总计 , 2022-Q2 ]
# After : sheets = [ 总计 , 2022-Q4 , 2022-Q2 ]
#   - "2022-Q4" reuses the worksheet that used to be named "2022-Q2"
#     (same sheetId/position) and gets the NEW fund-holdings data.
#   - A brand new worksheet named "2022-Q2" is appended right after it,
#     and receives the data that used to live on the old "2022-Q2" sheet.
#   - The "总计" (summary) sheet gets a new row 2 for 2022-Q4 and the
#     previous row 2 (2022-Q2) is pushed down to row 3.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the "总计" (totals) sheet
# ------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing 2022-Q2 summary row down to row 3 first so we don't
# clobber it while it is still in row 2.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(3, 2).Value = "2022-Q2"
$total.Cells.Item(3, 3).Value = 3
$total.Cells.Item(3, 4).Value = 0.17
$total.Cells.Item(3, 1).Font.Bold = $true
$total.Cells.Item(3, 1).Borders.LineStyle = 1
$total.Cells.Item(3, 1).HorizontalAlignment = -4108
$total.Cells.Item(3, 1).VerticalAlignment = -4160

# Write the new 2022-Q4 summary into row 2.
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 7
$total.Cells.Item(2, 4).Value = 0.45

# ------------------------------------------------------------------
# 2. Turn the existing "2022-Q2" detail sheet into the new "2022-Q4"
#    detail sheet (keeps its identity/position, rId2/sheet2.xml).
# ------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2022-Q2")
$q4.Name = "2022-Q4"
$q4.UsedRange.Clear()

$q4headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q4headers.Count; $c++) {
    $cell = $q4.Cells.Item(1, $c + 2)
    $cell.Value = $q4headers[$c]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$q4rows = @(
    @(0, "519020", "国泰金泰灵活配置混合A", "2.03", "93.85", "6.21", "0.1261", 7),
    @(1, "005433", "申万菱信医药先锋股票A", "2.12", "94.43", "5.42", "0.1149", 6),
    @(2, "014313", "鹏华创新增长一年持有期混合A", "3.47", "60.62", "3.12", "0.1083", 9),
    @(3, "005970", "国泰消费优选股票", "0.98", "93.52", "6.37", "0.0624", 4),
    @(4, "519022", "国泰金泰灵活配置混合C", "0.51", "93.85", "6.21", "0.0317", 7),
    @(5, "014314", "鹏华创新增长一年持有期混合C", "0.15", "60.62", "3.12", "0.0047", 9),
    @(6, "015171", "申万菱信医药先锋股票C", "0.01", "94.43", "5.42", "0.0005", 6)
)

for ($i = 0; $i -lt $q4rows.Count; $i++) {
    $row = $q4rows[$i]
    $r = $i + 2

    $aCell = $q4.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $q4.Cells.Item($r, 2).Value = "'" + $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = "'" + $row[3]
    $q4.Cells.Item($r, 5).Value = "'" + $row[4]
    $q4.Cells.Item($r, 6).Value = "'" + $row[5]
    $q4.Cells.Item($r, 7).Value = "'" + $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
}

# ------------------------------------------------------------------
# 3. Create the new "2022-Q2" detail sheet (after "2022-Q4") that
#    carries the data which previously lived on the old sheet2.
# ------------------------------------------------------------------
$q2 = $wb.Worksheets.Add($null, $q4)
$q2.Name = "2022-Q2"

$q2headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q2headers.Count; $c++) {
    $cell = $q2.Cells.Item(1, $c + 2)
    $cell.Value = $q2headers[$c]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$q2rows = @(
    @(0, "519020", "国泰金泰灵活配置混合A", "1.91", "93.56", "6.11", "0.1167", 9),
    @(1, "005970", "国泰消费优选股票", "0.90", "93.45", "6.05", "0.0544", 2),
    @(2, "519022", "国泰金泰灵活配置混合C", "0.06", "93.56", "6.11", "0.0037", 9)
)

for ($i = 0; $i -lt $q2rows.Count; $i++) {
    $row = $q2rows[$i]
    $r = $i + 2

    $aCell = $q2.Cells.Item($r, 1)
    $aCell.Value = $row[0]
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    $q2.Cells.Item($r, 2).Value = "'" + $row[1]
    $q2.Cells.Item($r, 3).Value = $row[2]
    $q2.Cells.Item($r, 4).Value = "'" + $row[3]
    $q2.Cells.Item($r, 5).Value = "'" + $row[4]
    $q2.Cells.Item($r, 6).Value = "'" + $row[5]
    $q2.Cells.Item($r, 7).Value = "'" + $row[6]
    $q2.Cells.Item($r, 8).Value = $row[7]
}
